$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.591135333333333
$ws.Range("H2").Value = 4.773406
$ws.Range("I2").Value = 0.4707829234247397
$ws.Range("J2").Value = 0.4707829234247397
$ws.Range("M2").Value = 0.006530999999999999
$ws.Range("N2").Value = 0.019593
$ws.Range("Q2").Value = 0.010391704862
$ws.Range("R2").Value = 0.09352534375799999
$ws.Range("S2").Value = 0.4707829234247397
$ws.Range("T2").Value = 0.4707829234247397

$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("I3").Value = 0.3035973020998604
$ws.Range("J3").Value = 0.3035973020998604
$ws.Range("M3").Value = 0.006530999999999999
$ws.Range("N3").Value = 0.019593
$ws.Range("Q3").Value = 0.006701376374
$ws.Range("R3").Value = 0.060312387366
$ws.Range("S3").Value = 0.3035973020998604
$ws.Range("T3").Value = 0.3035973020998604

$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 0.730693
$ws.Range("H4").Value = 2.192079
$ws.Range("I4").Value = 0.2161964349979826
$ws.Range("J4").Value = 0.2161964349979826
$ws.Range("M4").Value = 0.006530999999999999
$ws.Range("N4").Value = 0.019593
$ws.Range("Q4").Value = 0.004772155983
$ws.Range("R4").Value = 0.042949403847
$ws.Range("S4").Value = 0.2161964349979826
$ws.Range("T4").Value = 0.2161964349979826

$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.03184866666666667
$ws.Range("H5").Value = 0.095546
$ws.Range("I5").Value = 0.009423339477417213
$ws.Range("J5").Value = 0.009423339477417213
$ws.Range("M5").Value = 0.006530999999999999
$ws.Range("N5").Value = 0.019593
$ws.Range("Q5").Value = 0.000208003642
$ws.Range("R5").Value = 0.001872032778
$ws.Range("S5").Value = 0.009423339477417213
$ws.Range("T5").Value = 0.009423339477417213
